$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) from column AC into the new AD:AF columns so the
# appended quarters (31/12/2023, 31/03/2024, 30/06/2024) match the look of the
# existing header/data/blank-separator cells.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null
$ws.Range("AC2:AC80").Copy() | Out-Null
$ws.Range("AD2:AD80").PasteSpecial(-4122) | Out-Null
$ws.Range("AE2:AE80").PasteSpecial(-4122) | Out-Null
$ws.Range("AF2:AF80").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# New quarter headers
$ws.Range("AD1").Value = "31/12/2023"
$ws.Range("AE1").Value = "31/03/2024"
$ws.Range("AF1").Value = "30/06/2024"

# New quarter data (rows 2-80, mirrors the financial statement layout already
# present for the earlier quarters). Rows 57, 58, 71-73, 77-78 are the blank
# section-separator rows and are intentionally left empty (format-only, copied above).
$ws.Range("AD2").Value = 43481001.984
$ws.Range("AE2").Value = 43157999.616
$ws.Range("AF2").Value = 44302999.552
$ws.Range("AD3").Value = 23598999.552
$ws.Range("AE3").Value = 23079000.064
$ws.Range("AF3").Value = 24764000.256
$ws.Range("AD4").Value = 6665999.872
$ws.Range("AE4").Value = 5475999.744
$ws.Range("AF4").Value = 8224999.936
$ws.Range("AD5").Value = 0
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 0
$ws.Range("AD6").Value = 6135000.064
$ws.Range("AE6").Value = 6042999.808
$ws.Range("AF6").Value = 5877000.192
$ws.Range("AD7").Value = 5953999.872
$ws.Range("AE7").Value = 7107999.744
$ws.Range("AF7").Value = 6264000
$ws.Range("AD8").Value = 0
$ws.Range("AE8").Value = 0
$ws.Range("AF8").Value = 0
$ws.Range("AD9").Value = 3641999.872
$ws.Range("AE9").Value = 3152000
$ws.Range("AF9").Value = 3081999.872
$ws.Range("AD10").Value = 106000
$ws.Range("AE10").Value = 109000
$ws.Range("AF10").Value = 130000
$ws.Range("AD11").Value = 1096000
$ws.Range("AE11").Value = 1191000.064
$ws.Range("AF11").Value = 1186000
$ws.Range("AD12").Value = 7327000.064
$ws.Range("AE12").Value = 7322999.808
$ws.Range("AF12").Value = 7045000.192
$ws.Range("AD13").Value = 0
$ws.Range("AE13").Value = 0
$ws.Range("AF13").Value = 0
$ws.Range("AD14").Value = 0
$ws.Range("AE14").Value = 0
$ws.Range("AF14").Value = 0
$ws.Range("AD15").Value = 0
$ws.Range("AE15").Value = 0
$ws.Range("AF15").Value = 0
$ws.Range("AD16").Value = 391000
$ws.Range("AE16").Value = 364000
$ws.Range("AF16").Value = 373000
$ws.Range("AD17").Value = 0
$ws.Range("AE17").Value = 0
$ws.Range("AF17").Value = 0
$ws.Range("AD18").Value = 0
$ws.Range("AE18").Value = 0
$ws.Range("AF18").Value = 0
$ws.Range("AD19").Value = 4148999.936
$ws.Range("AE19").Value = 4264000
$ws.Range("AF19").Value = 3992999.936
$ws.Range("AD20").Value = 62000
$ws.Range("AE20").Value = 58000
$ws.Range("AF20").Value = 54000
$ws.Range("AD21").Value = 0
$ws.Range("AE21").Value = 0
$ws.Range("AF21").Value = 0
$ws.Range("AD22").Value = 4489999.872
$ws.Range("AE22").Value = 4440000
$ws.Range("AF22").Value = 4524000.256
$ws.Range("AD23").Value = 6953999.872
$ws.Range("AE23").Value = 6931999.744
$ws.Range("AF23").Value = 6749000.192
$ws.Range("AD24").Value = 1111000.064
$ws.Range("AE24").Value = 1384000
$ws.Range("AF24").Value = 1220999.936
$ws.Range("AD25").Value = 0
$ws.Range("AE25").Value = 0
$ws.Range("AF25").Value = 0
$ws.Range("AD26").Value = 43481001.984
$ws.Range("AE26").Value = 43157999.616
$ws.Range("AF26").Value = 44302999.552
$ws.Range("AD27").Value = 9996000.255999999
$ws.Range("AE27").Value = 10490000.384
$ws.Range("AF27").Value = 9952000
$ws.Range("AD28").Value = 302000
$ws.Range("AE28").Value = 216000
$ws.Range("AF28").Value = 245000
$ws.Range("AD29").Value = 4496000
$ws.Range("AE29").Value = 4353999.872
$ws.Range("AF29").Value = 3104000
$ws.Range("AD30").Value = 1034000
$ws.Range("AE30").Value = 421000
$ws.Range("AF30").Value = 27000
$ws.Range("AD31").Value = 1470000
$ws.Range("AE31").Value = 2951000.064
$ws.Range("AF31").Value = 3984000
$ws.Range("AD32").Value = 0
$ws.Range("AE32").Value = 0
$ws.Range("AF32").Value = 0
$ws.Range("AD33").Value = 1124000
$ws.Range("AE33").Value = 683000
$ws.Range("AF33").Value = 1154000
$ws.Range("AD34").Value = 1415000.064
$ws.Range("AE34").Value = 1709999.936
$ws.Range("AF34").Value = 1310999.936
$ws.Range("AD35").Value = 155000
$ws.Range("AE35").Value = 155000
$ws.Range("AF35").Value = 127000
$ws.Range("AD36").Value = 0
$ws.Range("AE36").Value = 0
$ws.Range("AF36").Value = 0
$ws.Range("AD37").Value = 17753999.36
$ws.Range("AE37").Value = 16112000
$ws.Range("AF37").Value = 17813000.192
$ws.Range("AD38").Value = 14048000
$ws.Range("AE38").Value = 12696000.512
$ws.Range("AF38").Value = 14721999.872
$ws.Range("AD39").Value = 0
$ws.Range("AE39").Value = 0
$ws.Range("AF39").Value = 0
$ws.Range("AD40").Value = 0
$ws.Range("AE40").Value = 0
$ws.Range("AF40").Value = 0
$ws.Range("AD41").Value = 0
$ws.Range("AE41").Value = 0
$ws.Range("AF41").Value = 0
$ws.Range("AD42").Value = 0
$ws.Range("AE42").Value = 0
$ws.Range("AF42").Value = 0
$ws.Range("AD43").Value = 3705999.872
$ws.Range("AE43").Value = 3416000
$ws.Range("AF43").Value = 3091000.064
$ws.Range("AD44").Value = 0
$ws.Range("AE44").Value = 0
$ws.Range("AF44").Value = 0
$ws.Range("AD45").Value = 0
$ws.Range("AE45").Value = 0
$ws.Range("AF45").Value = 0
$ws.Range("AD46").Value = 0
$ws.Range("AE46").Value = 0
$ws.Range("AF46").Value = 0
$ws.Range("AD47").Value = 15731000.32
$ws.Range("AE47").Value = 16556000.256
$ws.Range("AF47").Value = 16538000.384
$ws.Range("AD48").Value = 7578999.808
$ws.Range("AE48").Value = 7578999.808
$ws.Range("AF48").Value = 10033999.872
$ws.Range("AD49").Value = -1091000.064
$ws.Range("AE49").Value = -1063000
$ws.Range("AF49").Value = 2000
$ws.Range("AD50").Value = 0
$ws.Range("AE50").Value = 0
$ws.Range("AF50").Value = 0
$ws.Range("AD51").Value = 10632999.936
$ws.Range("AE51").Value = 11422000.128
$ws.Range("AF51").Value = 7850999.808
$ws.Range("AD52").Value = 0
$ws.Range("AE52").Value = 0
$ws.Range("AF52").Value = 0
$ws.Range("AD53").Value = -1390000
$ws.Range("AE53").Value = -1382000
$ws.Range("AF53").Value = -1348999.936
$ws.Range("AD54").Value = 0
$ws.Range("AE54").Value = 0
$ws.Range("AF54").Value = 0
$ws.Range("AD55").Value = 0
$ws.Range("AE55").Value = 0
$ws.Range("AF55").Value = 0
$ws.Range("AD56").Value = 0
$ws.Range("AE56").Value = 0
$ws.Range("AF56").Value = 0
$ws.Range("AD59").Value = 43663003.648
$ws.Range("AE59").Value = 39599001.6
$ws.Range("AF59").Value = 42109001.728
$ws.Range("AD60").Value = -41258004.48
$ws.Range("AE60").Value = -37488001.024
$ws.Range("AF60").Value = -40097001.472
$ws.Range("AD61").Value = 2404999.68
$ws.Range("AE61").Value = 2111000.064
$ws.Range("AF61").Value = 2012000
$ws.Range("AD62").Value = -674999.936
$ws.Range("AE62").Value = -666000
$ws.Range("AF62").Value = -644000
$ws.Range("AD63").Value = -232000
$ws.Range("AE63").Value = -224000
$ws.Range("AF63").Value = -238000
$ws.Range("AD64").Value = 0
$ws.Range("AE64").Value = 0
$ws.Range("AF64").Value = 0
$ws.Range("AD65").Value = 0
$ws.Range("AE65").Value = 0
$ws.Range("AF65").Value = 0
$ws.Range("AD66").Value = 2970000.128
$ws.Range("AE66").Value = 408000
$ws.Range("AF66").Value = 12000
$ws.Range("AD67").Value = 41000
$ws.Range("AE67").Value = -71000
$ws.Range("AF67").Value = 79000
$ws.Range("AD68").Value = -50000
$ws.Range("AE68").Value = -334000
$ws.Range("AF68").Value = -213000
$ws.Range("AD69").Value = 249000
$ws.Range("AE69").Value = 189000
$ws.Range("AF69").Value = 304000
$ws.Range("AD70").Value = -299000
$ws.Range("AE70").Value = -523000
$ws.Range("AF70").Value = -517000
$ws.Range("AD74").Value = 4458999.808
$ws.Range("AE74").Value = 1224000
$ws.Range("AF74").Value = 1008000
$ws.Range("AD75").Value = -1052999.936
$ws.Range("AE75").Value = -451000
$ws.Range("AF75").Value = -47000
$ws.Range("AD76").Value = -109000
$ws.Range("AE76").Value = 16000
$ws.Range("AF76").Value = -94000
$ws.Range("AD79").Value = 0
$ws.Range("AE79").Value = 0
$ws.Range("AF79").Value = 0
$ws.Range("AD80").Value = 3296999.936
$ws.Range("AE80").Value = 789000
$ws.Range("AF80").Value = 867000

Write-Host "VBBR3 quarters appended (AD:AF)"
